# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (column G), recalculated per-row, replacing the previous
# "Strike#" derived figures.
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 4
    9  = 2
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 2
    24 = 3
    25 = 1
    26 = 2
    27 = 2
    28 = 1
    29 = 0
    30 = 3
    31 = 4
    32 = 1
    33 = 2
    34 = 2
    35 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
